$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.660.53"
$ws.Range("E2").Value = "  +0.63%  "
$ws.Range("D3").Value = "1.643.78"
$ws.Range("E3").Value = "  +1.01%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'215.10"
$ws.Range("E5").Value = "  +1.07%  "
$ws.Range("E6").Value = "  +1.74%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  +0.99%  "
$ws.Range("D9").Value = "'0.0627"
$ws.Range("E9").Value = "  +0.99%  "
$ws.Range("D10").Value = "'19.25"
$ws.Range("E10").Value = "  +1.18%  "
$ws.Range("E11").Value = "  +0.02%  "
$ws.Range("D12").Value = "1.873.08"
$ws.Range("E12").Value = "  +1.05%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.664.12"
$ws.Range("E13").Value = "  +2.84%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'4.19"
$ws.Range("E14").Value = "  +2.14%  "
$ws.Range("E15").Value = "  +1.96%  "
$ws.Range("D16").Value = "'65.09"
$ws.Range("E16").Value = "  +2.61%  "
$ws.Range("D17").Value = "26.669.06"
$ws.Range("E17").Value = "  +0.71%  "
$ws.Range("D18").Value = "0.0₃0745"
$ws.Range("E18").Value = "  +0.75%  "
$ws.Range("D19").Value = "'216.36"
$ws.Range("E19").Value = "  +0.81%  "
$ws.Range("E20").Value = "  -0.04%  "
$ws.Range("D21").Value = "'4.36"
$ws.Range("E21").Value = "  +1.52%  "
$ws.Range("E22").Value = "  +2.00%  "
$ws.Range("E23").Value = "  +1.77%  "
$ws.Range("E24").Value = "  +15.31%  "
$ws.Range("D25").Value = "'145.72"
$ws.Range("E25").Value = "  -1.40%  "
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("E27").Value = "  -0.05%  "
$ws.Range("D28").Value = "'7.17"
$ws.Range("E28").Value = "  +4.63%  "
$ws.Range("D29").Value = "'15.71"
$ws.Range("E29").Value = "  +1.33%  "
$ws.Range("E30").Value = "  +2.66%  "
$ws.Range("D31").Value = "'1.18"
$ws.Range("E31").Value = "  +1.13%  "
$ws.Range("E32").Value = "  +2.92%  "
$ws.Range("D33").Value = "'3.05"
$ws.Range("E33").Value = "  +3.39%  "
$ws.Range("D34").Value = "1.275.72"
$ws.Range("E34").Value = "  +4.87%  "
$ws.Range("E35").Value = "  +2.63%  "
$ws.Range("E36").Value = "  +0.80%  "
$ws.Range("D37").Value = "'0.0179"
$ws.Range("E37").Value = "  +3.92%  "
$ws.Range("E38").Value = "  +6.89%  "
$ws.Range("E39").Value = "  +4.81%  "
$ws.Range("E40").Value = "  -0.08%  "
$ws.Range("D41").Value = "'0.816"
$ws.Range("E41").Value = "  +2.82%  "
$ws.Range("E42").Value = "  -0.15%  "
$ws.Range("D43").Value = "'5.46"
$ws.Range("E43").Value = "  +2.27%  "
$ws.Range("D44").Value = "1.783.05"
$ws.Range("E44").Value = "  +1.16%  "
$ws.Range("D45").Value = "'91.27"
$ws.Range("E45").Value = "  -0.87%  "
$ws.Range("D46").Value = "'59.72"
$ws.Range("E46").Value = "  +8.85%  "
$ws.Range("E47").Value = "  +2.71%  "
$ws.Range("E48").Value = "  -0.16%  "
$ws.Range("D49").Value = "'0.0516"
$ws.Range("E49").Value = "  +0.81%  "
$ws.Range("E50").Value = "  +2.69%  "
$ws.Range("D51").Value = "'0.0969"
$ws.Range("E51").Value = "  +3.00%  "
